$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of translation data (row 40)
$ws.Range("A40").Value = "NOT-VALID"
$ws.Range("B40").Value = "AT LEAST 10 COINS"
$ws.Range("C40").Value = "BẠN PHẢI NHẬP ÍT NHẤT 10 XU"

# Update selection / view to match the new active cell
$ws.Range("C40").Select()
